$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Prepend letters to the supplier names in the shared strings (A2:A26 contain these values)
$ws.Range("A2:A6").Value = "A Mono Packaging Materials"
$ws.Range("A7:A11").Value = "B Trio PET PLC"
$ws.Range("A12:A16").Value = "C Miami Oranges"
$ws.Range("A17:A21").Value = "D NO8DO Mango"
$ws.Range("A22:A26").Value = "E Seitan Vitamins"

$ws.Range("A22").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2:A6").Select()
